# Update faculty and major data in the worksheet.
# Adds 13 new faculty/program rows (rows 3-15) below the existing
# "GENERAL" row, following the same code/name/active layout as the
# existing data (columns A/B/C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - name entered before code (matches source order of entry)
$ws.Range("B3").Value = "Fakultas Ilmu Komputer dan Teknologi Informasi"
$ws.Range("A3").Value = "FIKTI"
$ws.Range("C3").Value = $true

# Row 4
$ws.Range("A4").Value = "FIKF"
$ws.Range("B4").Value = "Fakultas Ilmu Kesehatan dan Farmasi"
$ws.Range("C4").Value = $true

# Row 5
$ws.Range("A5").Value = "FK"
$ws.Range("B5").Value = "Fakultas Kedokteran"
$ws.Range("C5").Value = $true

# Row 6
$ws.Range("A6").Value = "FTSP"
$ws.Range("B6").Value = "Fakultas Teknik Sipil dan Perencanaan"
$ws.Range("C6").Value = $true

# Row 7
$ws.Range("A7").Value = "FE"
$ws.Range("B7").Value = "Fakultas Ekonomi"
$ws.Range("C7").Value = $true

# Row 8
$ws.Range("A8").Value = "FPSI"
$ws.Range("B8").Value = "Fakultas Psikologi"
$ws.Range("C8").Value = $true

# Row 9
$ws.Range("A9").Value = "FTI"
$ws.Range("B9").Value = "Fakultas Teknologi Industri"
$ws.Range("C9").Value = $true

# Row 10
$ws.Range("A10").Value = "FSB"
$ws.Range("B10").Value = "Fakultas Sastra dan Budaya"
$ws.Range("C10").Value = $true

# Row 11
$ws.Range("A11").Value = "FIKOM"
$ws.Range("B11").Value = "Fakultas Ilmu Komunikasi"
$ws.Range("C11").Value = $true

# Row 12
$ws.Range("A12").Value = "FTI_DIPLOMA"
$ws.Range("B12").Value = "Fakultas Teknologi Informasi"
$ws.Range("C12").Value = $true

# Row 13
$ws.Range("A13").Value = "FBK"
$ws.Range("B13").Value = "Fakultas Bisnis dan Kewirausahaan"
$ws.Range("C13").Value = $true

# Row 14
$ws.Range("A14").Value = "MAGISTER"
$ws.Range("B14").Value = "Program Magister"
$ws.Range("C14").Value = $true

# Row 15
$ws.Range("A15").Value = "PROFESI"
$ws.Range("B15").Value = "Program Profesi"
$ws.Range("C15").Value = $true

# Match the final selected cell left behind in the saved workbook.
$ws.Range("E7").Select()
